# The diff adds a second data row (row 2) to the "Đơn sale phụ" sheet,
# extending the used range from A1:T1 to A1:T2. The new row's numeric
# columns (B, I, K, L, M, N, O, P) hold 0; the remaining (text) columns
# are left blank/empty, matching the target row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Đơn sale phụ")

$ws.Cells.Item(2, 2).Value = 0   # B2
$ws.Cells.Item(2, 9).Value = 0   # I2
$ws.Cells.Item(2, 11).Value = 0  # K2
$ws.Cells.Item(2, 12).Value = 0  # L2
$ws.Cells.Item(2, 13).Value = 0  # M2
$ws.Cells.Item(2, 14).Value = 0  # N2
$ws.Cells.Item(2, 15).Value = 0  # O2
$ws.Cells.Item(2, 16).Value = 0  # P2
